$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.644.39"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.153.12"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.32"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.49"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +14.84%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.439"
$ws.Range("E10").Value = "  +5.41%  "
$ws.Range("E11").Value = "  +3.06%  "
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.696.11"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.82"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("E15").Value = "  +3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.684.02"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.26"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.156.74"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.14"
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.29"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.79"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.53"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  +12.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0865"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.28"
$ws.Range("E30").Value = "  +3.44%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.44"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.27"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.95"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.639.80"
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").Value = "  +5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.723"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.03"
$ws.Range("E44").Value = "  +3.55%  "
$ws.Range("E45").Value = "  +6.38%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.192.00"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.104"
$ws.Range("E48").Value = "  +13.49%  "
$ws.Range("E49").Value = "  +1.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.977"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  +0.00%  "
